$wb = $excel.ActiveWorkbook

# --- Items sheet: add a new row for "nonExistant Item" with quantity 2 ---
$wsItems = $wb.Worksheets.Item("Items")

# Bring the new row's formatting in line with the row above it (style index
# "1" used throughout the Items data rows), then fill in the values.
$wsItems.Range("A45:B45").Copy()
$wsItems.Range("A46:B46").PasteSpecial(-4122)
$wsItems.Range("A46").Value = "nonExistant Item"
$wsItems.Range("B46").Value = 2

# --- ClientItems sheet: add a new row ClientId=1, ItemId=45, Quantity=3 ---
$wsClientItems = $wb.Worksheets.Item("ClientItems")
$wsClientItems.Range("A17").Value = 1
$wsClientItems.Range("B17").Value = 45
$wsClientItems.Range("C17").Value = 3

# --- Leave the view/selection state matching what was recorded after the edit ---
# (Items keeps its own remembered selection even while ClientItems stays the
# active/visible tab, exactly like Excel preserves per-sheet selection.)
$wsItems.Activate()
$wsItems.Range("A46").Select()

$wsClientItems.Activate()
$wsClientItems.Range("F16:F17").Select()
